$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 15: Wim Leerman (written first so "Wim"/"Leerman" get the lowest new shared-string ids) ---
$ws.Range("A15").Value = "Wim"
$ws.Range("B15").Value = "Leerman"
$ws.Range("C15").Value = 1997
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = "Ad-Senator"

# --- Row 6: Maarten Verpoest ---
$ws.Range("A6").Value = "Maarten"
$ws.Range("B6").Value = "Verpoest"
$ws.Range("C6").Value = 2010
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = "Ad-Senator"

# --- Row 7: Nicolas Van Vlasselaer ---
$ws.Range("A7").Value = "Nicolas"
$ws.Range("B7").Value = "Van Vlasselaer"
$ws.Range("C7").Value = 2013
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = "Ad-Senator"

# --- Row 11: Arnaud Jacquet ---
$ws.Range("A11").Value = "Arnaud"
$ws.Range("B11").Value = "Jacquet"
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = "Actief lid"

# --- Row 14: Els Denorme ---
$ws.Range("A14").Value = "Els"
$ws.Range("B14").Value = "Denorme"
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = "Senator"

# --- Row 16: Maarten Wuijts ---
$ws.Range("A16").Value = "Maarten"
$ws.Range("B16").Value = "Wuijts"
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = "Actief lid"

# --- Row 13: Philip Vermeylen ---
$ws.Range("A13").Value = "Philip"
$ws.Range("B13").Value = "Vermeylen"
$ws.Range("C13").Value = 2010
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = "Ad-Senator"

# --- Row 12: Tom Van Iseghem ---
$ws.Range("A12").Value = "Tom"
$ws.Range("B12").Value = "Van Iseghem"
$ws.Range("C12").Value = 2012
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = "Ad-Senator"

# --- Row 10: fill in previously-blank F10 ---
$ws.Range("F10").Value = 0

# --- E4/E5: clear the stored 0 but keep the numeric style ---
$ws.Range("E4").ClearContents()
$ws.Range("E5").ClearContents()

# --- Selection moves to E22 ---
$ws.Range("E22").Select()
